$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 90 (shifts existing rows 90-94 down to 91-95)
$ws.Rows.Item(90).Insert()

# Populate the new row 90 with data (copy of surrounding row pattern, new values)
$ws.Cells.Item(90, 1).Value = 10
$ws.Cells.Item(90, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(90, 3).Value = "La Araucanía"
$ws.Cells.Item(90, 4).Value = 44746
$ws.Cells.Item(90, 5).Value = 9
$ws.Cells.Item(90, 6).Value = 100114002
$ws.Cells.Item(90, 7).Value = "Camote"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 30
$ws.Cells.Item(90, 11).Value = 15000
$ws.Cells.Item(90, 12).Value = 15000
$ws.Cells.Item(90, 13).Value = 15000
$ws.Cells.Item(90, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(90, 15).Value = "Perú"
$ws.Cells.Item(90, 16).Value = 750
$ws.Cells.Item(90, 17).Value = 20
$ws.Cells.Item(90, 18).Value = "Hortaliza"

# Apply the date number format (matching the style used by other date cells in column D)
$ws.Cells.Item(90, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
